$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header cell formatting (style) onto the two new header cells, then set their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 8),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 3),
    @(1, 7),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(5, 7),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 9),
    @(2, 5),
    @(5, 8),
    @(7, 7),
    @(5, 8),
    @(4, 7),
    @(7, 8),
    @(1, 3),
    @(1, 3),
    @(7, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
